$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1 and 2 (Recursor HumanosPT group: remove K'IAM and TalentEd entries)
$ws.Range("A1:D1").ClearContents()
$ws.Range("A1:D1").VerticalAlignment = -4108
$ws.Range("A2:D2").ClearContents()
$ws.Range("A2:D2").VerticalAlignment = -4108

# Row 3 keeps Docubase entry but now carries the "Liens :" label in column B
$ws.Range("B3").Value = "Liens :"
$ws.Range("B3").VerticalAlignment = -4108

# Rows 13 and 14 (OutilsPT group: remove Easy-It and K'IAM entries)
$ws.Range("A13:D13").ClearContents()
$ws.Range("A13:D13").VerticalAlignment = -4108
$ws.Range("A14:D14").ClearContents()
$ws.Range("A14:D14").VerticalAlignment = -4108

# Row 15 keeps Forecast & Replenishment entry but now carries the "Liens :" label
$ws.Range("B15").Value = "Liens :"
$ws.Range("B15").VerticalAlignment = -4108

# Row 16 (Programa de DecisionPT group: remove Kiperf entry)
$ws.Range("A16:D16").ClearContents()
$ws.Range("A16:D16").VerticalAlignment = -4108

# Row 17 keeps Zone de lancement BI entry but now carries the "Liens :" label
$ws.Range("B17").Value = "Liens :"
$ws.Range("B17").VerticalAlignment = -4108
